# Actualización automática 2025-08-25 14:40:08
#
# Inserts a new client row ("LATACELA ZUÑIGA JUAN FERNANDO") for asesor
# "GUERRERO FAREZ FABIAN MAURICIO" at row 26 (alphabetically before
# "MADECOR-HOME CENTER S.A.S.") on both the "VENTAS POR GRUPO" and
# "VENTA MENSUAL" sheets, pushing the existing rows 26-54 down to 27-55.
# Also refreshes the trailing summary row (counts "X de 52" -> "X de 53"
# on sheet 1, grand totals on sheet 2).

$wb = $excel.ActiveWorkbook

$asesor = "GUERRERO FAREZ FABIAN MAURICIO"
$cliente = "LATACELA ZUÑIGA JUAN FERNANDO"

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"  (columns A:R)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a blank row above the current row 26 (MADECOR...), shifting
# everything below (including the trailing summary row) down by one.
$ws1.Rows.Item(26).Insert()

$ws1.Cells.Item(26, 1).Value = $asesor
$ws1.Cells.Item(26, 2).Value = $cliente

$row1Values = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row1Values.Length; $i++) {
    $ws1.Cells.Item(26, 3 + $i).Value2 = $row1Values[$i]
}

# Refresh the "N de 52" -> "N de 53" labels on the (now) last row, 55.
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(55, $col)
    $oldText = $cell.Value2
    $cell.Value = $oldText.Replace("de 52", "de 53")
}

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"  (columns A:G)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(26).Insert()

$ws2.Cells.Item(26, 1).Value = $asesor
$ws2.Cells.Item(26, 2).Value = $cliente
$ws2.Cells.Item(26, 3).Value2 = 171.19
$ws2.Cells.Item(26, 4).Value2 = 0
$ws2.Cells.Item(26, 5).Value2 = 0
$ws2.Cells.Item(26, 6).Value2 = 0
$ws2.Cells.Item(26, 7).Value2 = 1500

# Update the grand-total row (now row 55) to include the new row's amounts.
$ws2.Cells.Item(55, 3).Value2 = 97690.81999999999
$ws2.Cells.Item(55, 4).Value2 = 88077.29000000001
$ws2.Cells.Item(55, 5).Value2 = 93548.25
$ws2.Cells.Item(55, 6).Value2 = 62018.24
$ws2.Cells.Item(55, 7).Value2 = 130170.11
